$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (trade record updated by the fixed trading module)
$ws.Range("A2").Value = 42636.606770833336
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 19.29
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = $true

# Remove now-obsolete rows 3 and 4 entirely
$ws.Range("A3:G4").Delete()

# Restore the workbook window view size recorded by Excel on save
$excel.ActiveWindow.Width = 15345
$excel.ActiveWindow.Height = 6705
